$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)  # 展览
$ws4 = $wb.Worksheets.Item(4)  # 全部类型

# Sheet 1 (展览) updates
$ws1.Range("F2").Value = 288
$ws1.Range("F3").Value = 1464
$ws1.Range("G6").Value = 40
$ws1.Range("F7").Value = 105
$ws1.Range("F9").Value = 202
$ws1.Range("F10").Value = 152
$ws1.Range("F11").Value = 7
$ws1.Range("F12").Value = 4858
$ws1.Range("F14").Value = 7146
$ws1.Range("F17").Value = 95
$ws1.Range("F19").Value = 61
$ws1.Range("F21").Value = 4217
$ws1.Range("F22").Value = 1575
$ws1.Range("F24").Value = 81
$ws1.Range("F25").Value = 2791
$ws1.Range("F28").Value = 185
$ws1.Range("F29").Value = 414
$ws1.Range("F30").Value = 403
$ws1.Range("F31").Value = 422
$ws1.Range("F32").Value = 259
$ws1.Range("F34").Value = 1660
$ws1.Range("F35").Value = 1102
$ws1.Range("F37").Value = 1174
$ws1.Range("F38").Value = 95
$ws1.Range("F39").Value = 561
$ws1.Range("F43").Value = 34
$ws1.Range("F44").Value = 97
$ws1.Range("F45").Value = 1886
$ws1.Range("F47").Value = 31

# Sheet 4 (全部类型) updates
$ws4.Range("F2").Value = 288
$ws4.Range("F3").Value = 1464
$ws4.Range("G6").Value = 40
$ws4.Range("F7").Value = 105
$ws4.Range("F9").Value = 202
$ws4.Range("F10").Value = 152
$ws4.Range("F11").Value = 7
$ws4.Range("F12").Value = 4858
$ws4.Range("F14").Value = 7146
$ws4.Range("F17").Value = 95
$ws4.Range("F19").Value = 61
$ws4.Range("F21").Value = 4217
$ws4.Range("F22").Value = 1575
$ws4.Range("F24").Value = 81
$ws4.Range("F25").Value = 2791
$ws4.Range("F28").Value = 185
$ws4.Range("F29").Value = 414
$ws4.Range("F30").Value = 403
$ws4.Range("F31").Value = 422
$ws4.Range("F32").Value = 259
$ws4.Range("F34").Value = 1660
$ws4.Range("F35").Value = 1102
$ws4.Range("F37").Value = 1174
$ws4.Range("F38").Value = 95
$ws4.Range("F39").Value = 561
$ws4.Range("F43").Value = 34
$ws4.Range("F44").Value = 97
$ws4.Range("F45").Value = 1887
$ws4.Range("F47").Value = 31
